$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5 and append new rows 6-10 for the added "Neutro" cluster ---

# Row 2: M1 | Npy | Npy1r | ECs
$ws.Cells.Item(2, 1).Value = "M1"
$ws.Cells.Item(2, 2).Value = "Npy"
$ws.Cells.Item(2, 3).Value = "Npy1r"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2.709067
$ws.Cells.Item(2, 8).Value = 8.127201
$ws.Cells.Item(2, 9).Value = 0.670840434061904
$ws.Cells.Item(2, 10).Value = 0.670840434061904
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.252158666666667
$ws.Cells.Item(2, 14).Value = 3.756476
$ws.Cells.Item(2, 15).Value = 0.1749215513854956
$ws.Cells.Item(2, 16).Value = 0.1749215513854956
$ws.Cells.Item(2, 17).Value = 3.392181722630666
$ws.Cells.Item(2, 18).Value = 30.529635503676
$ws.Cells.Item(2, 19).Value = 0.1173444494582275
$ws.Cells.Item(2, 20).Value = 0.1173444494582275

# Row 3: M1 | Npy | Npy1r | Neutro
$ws.Cells.Item(3, 1).Value = "M1"
$ws.Cells.Item(3, 2).Value = "Npy"
$ws.Cells.Item(3, 3).Value = "Npy1r"
$ws.Cells.Item(3, 4).Value = "Neutro"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.709067
$ws.Cells.Item(3, 8).Value = 8.127201
$ws.Cells.Item(3, 9).Value = 0.670840434061904
$ws.Cells.Item(3, 10).Value = 0.670840434061904
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.002279333333333333
$ws.Cells.Item(3, 14).Value = 0.006838
$ws.Cells.Item(3, 15).Value = 0.000318413738933516
$ws.Cells.Item(3, 16).Value = 0.0003184137389335161
$ws.Cells.Item(3, 17).Value = 0.006174866715333333
$ws.Cells.Item(3, 18).Value = 0.055573800438
$ws.Cells.Item(3, 19).Value = 0.0002136048108374337
$ws.Cells.Item(3, 20).Value = 0.0002136048108374337

# Row 4: M1 | Npy | Npy1r | sCs
$ws.Cells.Item(4, 1).Value = "M1"
$ws.Cells.Item(4, 2).Value = "Npy"
$ws.Cells.Item(4, 3).Value = "Npy1r"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2.709067
$ws.Cells.Item(4, 8).Value = 8.127201
$ws.Cells.Item(4, 9).Value = 0.670840434061904
$ws.Cells.Item(4, 10).Value = 0.670840434061904
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 5.903963333333333
$ws.Cells.Item(4, 14).Value = 17.71189
$ws.Cells.Item(4, 15).Value = 0.8247600348755709
$ws.Cells.Item(4, 16).Value = 0.8247600348755709
$ws.Cells.Item(4, 17).Value = 15.99423223554333
$ws.Cells.Item(4, 18).Value = 143.94809011989
$ws.Cells.Item(4, 19).Value = 0.553282379792839
$ws.Cells.Item(4, 20).Value = 0.553282379792839

# Row 5: M2 | Npy | Npy1r | ECs
$ws.Cells.Item(5, 1).Value = "M2"
$ws.Cells.Item(5, 2).Value = "Npy"
$ws.Cells.Item(5, 3).Value = "Npy1r"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.267349666666667
$ws.Cells.Item(5, 8).Value = 3.802049
$ws.Cells.Item(5, 9).Value = 0.3138310719132735
$ws.Cells.Item(5, 10).Value = 0.3138310719132735
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.252158666666667
$ws.Cells.Item(5, 14).Value = 3.756476
$ws.Cells.Item(5, 15).Value = 0.1749215513854956
$ws.Cells.Item(5, 16).Value = 0.1749215513854956
$ws.Cells.Item(5, 17).Value = 1.586922868813778
$ws.Cells.Item(5, 18).Value = 14.282305819324
$ws.Cells.Item(5, 19).Value = 0.05489581797204283
$ws.Cells.Item(5, 20).Value = 0.05489581797204283

# Row 6: M2 | Npy | Npy1r | Neutro
$ws.Cells.Item(6, 1).Value = "M2"
$ws.Cells.Item(6, 2).Value = "Npy"
$ws.Cells.Item(6, 3).Value = "Npy1r"
$ws.Cells.Item(6, 4).Value = "Neutro"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.267349666666667
$ws.Cells.Item(6, 8).Value = 3.802049
$ws.Cells.Item(6, 9).Value = 0.3138310719132735
$ws.Cells.Item(6, 10).Value = 0.3138310719132735
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.002279333333333333
$ws.Cells.Item(6, 14).Value = 0.006838
$ws.Cells.Item(6, 15).Value = 0.000318413738933516
$ws.Cells.Item(6, 16).Value = 0.0003184137389335161
$ws.Cells.Item(6, 17).Value = 0.002888712340222223
$ws.Cells.Item(6, 18).Value = 0.025998411062
$ws.Cells.Item(6, 19).Value = 0.00009992812500141857
$ws.Cells.Item(6, 20).Value = 0.00009992812500141858

# Row 7: M2 | Npy | Npy1r | sCs
$ws.Cells.Item(7, 1).Value = "M2"
$ws.Cells.Item(7, 2).Value = "Npy"
$ws.Cells.Item(7, 3).Value = "Npy1r"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.267349666666667
$ws.Cells.Item(7, 8).Value = 3.802049
$ws.Cells.Item(7, 9).Value = 0.3138310719132735
$ws.Cells.Item(7, 10).Value = 0.3138310719132735
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 5.903963333333333
$ws.Cells.Item(7, 14).Value = 17.71189
$ws.Cells.Item(7, 15).Value = 0.8247600348755709
$ws.Cells.Item(7, 16).Value = 0.8247600348755709
$ws.Cells.Item(7, 17).Value = 7.482385962512224
$ws.Cells.Item(7, 18).Value = 67.34147366261
$ws.Cells.Item(7, 19).Value = 0.2588353258162293
$ws.Cells.Item(7, 20).Value = 0.2588353258162293

# Row 8: Neutro | Npy | Npy1r | ECs
$ws.Cells.Item(8, 1).Value = "Neutro"
$ws.Cells.Item(8, 2).Value = "Npy"
$ws.Cells.Item(8, 3).Value = "Npy1r"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.06190133333333334
$ws.Cells.Item(8, 8).Value = 0.185704
$ws.Cells.Item(8, 9).Value = 0.01532849402482255
$ws.Cells.Item(8, 10).Value = 0.01532849402482255
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.252158666666667
$ws.Cells.Item(8, 14).Value = 3.756476
$ws.Cells.Item(8, 15).Value = 0.1749215513854956
$ws.Cells.Item(8, 16).Value = 0.1749215513854956
$ws.Cells.Item(8, 17).Value = 0.07751029101155556
$ws.Cells.Item(8, 18).Value = 0.6975926191040001
$ws.Cells.Item(8, 19).Value = 0.002681283955225259
$ws.Cells.Item(8, 20).Value = 0.002681283955225259

# Row 9: Neutro | Npy | Npy1r | Neutro
$ws.Cells.Item(9, 1).Value = "Neutro"
$ws.Cells.Item(9, 2).Value = "Npy"
$ws.Cells.Item(9, 3).Value = "Npy1r"
$ws.Cells.Item(9, 4).Value = "Neutro"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.06190133333333334
$ws.Cells.Item(9, 8).Value = 0.185704
$ws.Cells.Item(9, 9).Value = 0.01532849402482255
$ws.Cells.Item(9, 10).Value = 0.01532849402482255
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.002279333333333333
$ws.Cells.Item(9, 14).Value = 0.006838
$ws.Cells.Item(9, 15).Value = 0.000318413738933516
$ws.Cells.Item(9, 16).Value = 0.0003184137389335161
$ws.Cells.Item(9, 17).Value = 0.0001410937724444444
$ws.Cells.Item(9, 18).Value = 0.001269843952
$ws.Cells.Item(9, 19).Value = 0.000004880803094663807
$ws.Cells.Item(9, 20).Value = 0.000004880803094663808

# Row 10: Neutro | Npy | Npy1r | sCs
$ws.Cells.Item(10, 1).Value = "Neutro"
$ws.Cells.Item(10, 2).Value = "Npy"
$ws.Cells.Item(10, 3).Value = "Npy1r"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.06190133333333334
$ws.Cells.Item(10, 8).Value = 0.185704
$ws.Cells.Item(10, 9).Value = 0.01532849402482255
$ws.Cells.Item(10, 10).Value = 0.01532849402482255
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 5.903963333333333
$ws.Cells.Item(10, 14).Value = 17.71189
$ws.Cells.Item(10, 15).Value = 0.8247600348755709
$ws.Cells.Item(10, 16).Value = 0.8247600348755709
$ws.Cells.Item(10, 17).Value = 0.3654632022844445
$ws.Cells.Item(10, 18).Value = 3.28916882056
$ws.Cells.Item(10, 19).Value = 0.01264232926650263
$ws.Cells.Item(10, 20).Value = 0.01264232926650263

